$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Hoja1: a new article ("mouse") was parsed/added as row 4
$ws1.Range("A4").Value = "mouse"
$ws1.Range("B4").Value = "mse451"
$ws1.Range("C4").Value = 2000
$ws1.Range("D4").Value = 15
$ws1.Range("E4").Value = 3
$ws1.Range("F4").Value = 1

# Hoja2: re-parsed with new data - "monitor" and "pc" rows are gone,
# "tv" and "celular" move up to rows 2-3, and a new "teclado" row is added as row 4
$ws2.Range("A2").Value = "tv"
$ws2.Range("B2").Value = "tv54165"
$ws2.Range("C2").Value = 2500
$ws2.Range("D2").Value = 10
$ws2.Range("E2").Value = 3
$ws2.Range("F2").Value = 1

$ws2.Range("A3").Value = "celular"
$ws2.Range("B3").Value = "cell516"
$ws2.Range("C3").Value = 52000
$ws2.Range("D3").Value = 7
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 2

$ws2.Range("A4").Value = "teclado"
$ws2.Range("B4").Value = "tcl4651"
$ws2.Range("C4").Value = 2100
$ws2.Range("D4").Value = 9
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 3

$ws2.Range("A5:F5").ClearContents()
